$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the systematic-uncertainty header columns:
#   M1: syst0_c -> syst1_u
#   N1: syst1_c -> syst2_c
#   O1: syst2_c -> syst3_c
# Order matters for shared-string append order (new strings are appended in the
# order they are first introduced, once unreferenced strings are dropped):
# touch N1/O1 first (O1 introduces the new "syst3_c" string) then M1 last
# (introduces "syst1_u" last) to match the target shared-string table order.
$ws.Range("N1").Value = "syst2_c"
$ws.Range("O1").Value = "syst3_c"
$ws.Range("M1").Value = "syst1_u"

# Move the active selection from L17 to L12
[void]$ws.Range("L12").Select()
